$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns (rows 2-51) are treated as text so that
# values like "1.78" or "46.60" are not auto-converted to numbers,
# then clear the temporary formatting so cell styles remain unchanged.
$rngText = $ws.Range("D2:E51")
$rngText.NumberFormat = "@"

$ws.Range("D2").Value = '33.732.74'
$ws.Range("E2").Value = '  +8.66%  '
$ws.Range("D3").Value = '1.774.73'
$ws.Range("E3").Value = '  +4.88%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = '224.55'
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("D6").Value = '0.553'
$ws.Range("E6").Value = '  +3.61%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").Value = '30.55'
$ws.Range("E8").Value = '  +3.54%  '
$ws.Range("D9").Value = '46.60'
$ws.Range("E9").Value = '  +4.05%  '
$ws.Range("D10").Value = '0.276'
$ws.Range("E10").Value = '  +3.52%  '
$ws.Range("D11").Value = '0.0658'
$ws.Range("E11").Value = '  +2.80%  '
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("D13").Value = '2.035.28'
$ws.Range("E13").Value = '  +5.28%  '
$ws.Range("D14").Value = '1.781.21'
$ws.Range("E14").Value = '  +5.47%  '
$ws.Range("D15").Value = '0.623'
$ws.Range("E15").Value = '  +2.25%  '
$ws.Range("D16").Value = '33.775.34'
$ws.Range("E16").Value = '  +8.73%  '
$ws.Range("D17").Value = '9.90'
$ws.Range("E17").Value = '  -3.67%  '
$ws.Range("D18").Value = '4.16'
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("D19").Value = '68.21'
$ws.Range("E19").Value = '  +1.92%  '
$ws.Range("D20").Value = '250.67'
$ws.Range("E20").Value = '  +1.11%  '
$ws.Range("D21").Value = '0.0₃0735'
$ws.Range("E21").Value = '  +2.19%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").Value = '10.23'
$ws.Range("E23").Value = '  +1.77%  '
$ws.Range("D24").Value = '4.18'
$ws.Range("E24").Value = '  -3.10%  '
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").Value = '158.24'
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = '16.38'
$ws.Range("E27").Value = '  +2.76%  '
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("D29").Value = '6.89'
$ws.Range("E29").Value = '  +2.26%  '
$ws.Range("E30").Value = '  +0.21%  '
$ws.Range("D31").Value = '3.79'
$ws.Range("E31").Value = '  +6.88%  '
$ws.Range("D32").Value = '0.0511'
$ws.Range("E32").Value = '  +2.07%  '
$ws.Range("D33").Value = '1.18'
$ws.Range("E33").Value = '  +3.07%  '
$ws.Range("D34").Value = '3.53'
$ws.Range("E34").Value = '  +5.49%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '1.475.33'
$ws.Range("E35").Value = '  -2.55%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = '1.78'
$ws.Range("E36").Value = '  +2.09%  '
$ws.Range("E37").Value = '  +3.01%  '
$ws.Range("D38").Value = '0.629'
$ws.Range("E38").Value = '  +2.22%  '
$ws.Range("D39").Value = '82.91'
$ws.Range("E39").Value = '  -0.29%  '
$ws.Range("D40").Value = '0.0184'
$ws.Range("E40").Value = '  +2.27%  '
$ws.Range("D41").Value = '2.36'
$ws.Range("E41").Value = '  +2.66%  '
$ws.Range("D42").Value = '2.69'
$ws.Range("E42").Value = '  +0.91%  '
$ws.Range("D43").Value = '0.881'
$ws.Range("E43").Value = '  +4.17%  '
$ws.Range("D44").Value = '2.07'
$ws.Range("E44").Value = '  +1.58%  '
$ws.Range("D45").Value = '0.0506'
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("E46").Value = '  +3.22%  '
$ws.Range("D47").Value = '1.935.03'
$ws.Range("E47").Value = '  +6.16%  '
$ws.Range("D48").Value = '5.75'
$ws.Range("E48").Value = '  +3.17%  '
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("D50").Value = '11.89'
$ws.Range("E50").Value = '  +15.11%  '
$ws.Range("D51").Value = '50.38'
$ws.Range("E51").Value = '  -3.23%  '

# Clear the temporary number formatting applied above so the cell
# styles match the original (unstyled) state.
$rngText.ClearFormats()
